$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that is the same for
# every data row (2..244). The update bumps that date forward by one day,
# from serial 45178 (2023-09-09) to serial 45179 (2023-09-10).
$ws.Range("C2:C244").Value = 45179
